$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for all data rows (2..205)
# from serial date 45172 (2023-09-03) to 45175 (2023-09-06).
$ws.Range("C2:C205").Value = 45175
